$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.4159458605706507
$ws.Range("D2").Value = 0.6814791206859918

# Row 3
$ws.Range("C3").Value = -0.2488675716156888
$ws.Range("D3").Value = 0.8057718321380909

# Row 4
$ws.Range("C4").Value = 1.97825999529021
$ws.Range("D4").Value = 0.06055600410472128
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = -0.1496930492389824
$ws.Range("D5").Value = 0.8823706442362671

# Row 6
$ws.Range("C6").Value = 0.1805095025562502
$ws.Range("D6").Value = 0.8584050862060639

# Row 7
$ws.Range("C7").Value = 2.292397272759771
$ws.Range("D7").Value = 0.03181567353629799

# Row 8
$ws.Range("C8").Value = 0.2608706927948466
$ws.Range("D8").Value = 0.7966186243036313

# Row 9
$ws.Range("C9").Value = 2.008927248542865
$ws.Range("D9").Value = 0.05697249625604184

# Row 10
$ws.Range("C10").Value = 0.1331274138889595
$ws.Range("D10").Value = 0.8953028332992468

# Row 11
$ws.Range("C11").Value = -2.045563143145746
$ws.Range("D11").Value = 0.05293965659209321
$ws.Range("G11").Value = "No"
